# Vip-Vipr2.xlsx update: refresh the LR-pair table with newly computed TPM
# values. The previous self-referencing "ECs -> ECs" row is dropped and the
# remaining Target-cluster rows (FAPs / MuSCs / Resolving-Mac) get their
# expression / specificity metrics recalculated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old row 2 (Sending cluster ECs -> Target cluster ECs). Rows 3-5
# shift up to become rows 2-4, which already carry the correct Target
# cluster labels (FAPs, MuSCs, Resolving-Mac) for the refreshed data below.
$ws.Rows.Item(2).Delete()

# New TPM-derived values, keyed by (row, column).
$values = @{
    "E2" = 3;  "F2" = 1;  "G2" = 1.164924666666667;  "H2" = 3.494774
    "I2" = 1;  "J2" = 1;  "K2" = 3;  "L2" = 1
    "M2" = 6.804012333333334;   "N2" = 20.412037
    "O2" = 0.5269116569106099;  "P2" = 0.5269116569106099
    "Q2" = 7.926161799404222;   "R2" = 71.33545619463801
    "S2" = 0.5269116569106099;  "T2" = 0.5269116569106099

    "E3" = 3;  "F3" = 1;  "G3" = 1.164924666666667;  "H3" = 3.494774
    "I3" = 1;  "J3" = 1;  "K3" = 3;  "L3" = 1
    "M3" = 6.095937333333333;   "N3" = 18.287812
    "O3" = 0.4720773983600821;  "P3" = 0.472077398360082
    "Q3" = 7.101307766054221;   "R3" = 63.911769894488
    "S3" = 0.4720773983600821;  "T3" = 0.472077398360082

    "E4" = 3;  "F4" = 1;  "G4" = 1.164924666666667;  "H4" = 3.494774
    "I4" = 1;  "J4" = 1;  "K4" = 1;  "L4" = 0.3333333333333333
    "M4" = 0.01305433333333333; "N4" = 0.039163
    "O4" = 0.001010944729308016; "P4" = 0.001010944729308016
    "Q4" = 0.01520731490688889;  "R4" = 0.136865834162
    "S4" = 0.001010944729308016; "T4" = 0.001010944729308016
}

foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value2 = $values[$ref]
}
